# Adds security check to Email Adapter
# Update JsonDemo sheet: rename project/value labels, update G/H columns,
# and move the selection to H5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value updates (rows 2-5) ---
# Column B: "Project Love" -> "Excel Security"
$ws.Range("B2:B5").Value = "Excel Security"

# Column F: "Chez Martha" -> "Laura's Test Environment"
$ws.Range("F2:F5").Value = "Laura's Test Environment"

# Columns G/H: per-row relabeling, written row by row so new shared
# strings are interned in the same order Excel would (G then H, top to
# bottom). H2 "Done" stays the same.
$ws.Range("G2").Value = "Working Environment"
$ws.Range("G3").Value = "Date"
$ws.Range("H3").Value = "April 10 2016"
$ws.Range("G4").Value = "Test Number"
$ws.Range("H4").Value = 2
$ws.Range("G5").Value = "Test Iteration"
$ws.Range("H5").Value = "2B"

# --- Selection moves to H5 ---
$ws.Range("H5").Select()
